$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Trend_instructions")
$ws2.Activate()

# Delete columns G:J (L1_agg_fuel, L2_CEDS_fuel, L3_agg_sector, L4_CEDS_sector).
# This shifts the old K,L,M columns (override_normalization, start_continuity,
# end_continuity) left into G,H,I - carrying their column-width formatting
# along with them.
$ws2.Range("G1:J1").EntireColumn.Delete() | Out-Null

# Add the new header in column J: user_pct_breakdowns
$ws2.Cells.Item(1, 10).Value2 = "user_pct_breakdowns"

# New column I (end_continuity) gets a best-fit width in the target file.
$ws2.Columns.Item(9).ColumnWidth = 12.5

# Update the selection to the new layout's corresponding cell
$ws2.Range("B1").Select() | Out-Null
$ws2.Range("J8").Select() | Out-Null

$wb.Save()
